# Update the project monitoring dashboard from the live data source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "ACTUAL END" (column K) for every data row (2..127) moves from 45812 to 45816.
$ws.Range("K2:K127").Value = 45816

# 2. "% COMPLETE" (column L) updates for the rows whose progress changed.
$ws.Range("L52").Value = 0.4668
$ws.Range("L55").Value = 0.7036
$ws.Range("L58").Value = 1
$ws.Range("L62").Value = 0.9596
$ws.Range("L65").Value = 0.6233
$ws.Range("L69").Value = 0.4658
$ws.Range("L70").Value = 0.7808
$ws.Range("L71").Value = 0.1244
$ws.Range("L81").Value = 0.6317
$ws.Range("L84").Value = 0.4143
$ws.Range("L88").Value = 0.4425
$ws.Range("L115").Value = 0.2008

# 3. "PRIORITY" (column H) bumped to HIGH for a couple of tasks.
$ws.Range("H65").Value = "HIGH"
$ws.Range("H81").Value = "HIGH"

# 4. Update the active sheet view (frozen pane / current selection) to where the
#    user was last working.
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollColumn = 10
$appWin.Panes.Item(2).Activate()
$ws.Range("T24").Select()
